# The source workbook contains a small lookup/header table in A1:G5 where
# row 1 holds the column headers (B1:G1) and rows 2-5 hold data rows, with
# column A holding a styled (bordered/bold/centered) row index.
#
# This edit removes all of the data rows' contents: the B:G values are
# cleared out entirely (freeing up the now-unused shared strings that held
# the timestamps / company names / etc.), the index cells in column A keep
# their formatting but lose their numeric value, and the now fully blank
# 5th row is removed outright so the sheet's used range shrinks from
# A1:G5 down to A1:G4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out all the data-row contents (columns B:G) for rows 2-5.
$ws.Range("B2:G5").ClearContents()

# Column A keeps its styling (border/bold/center) on rows 2-4 but loses its
# value.
$ws.Range("A2:A4").ClearContents()

# The last row (5) is now entirely empty - delete it so the sheet shrinks
# to A1:G4, shifting nothing else below it.
$ws.Rows("5:5").Delete()
